# Weekly update: insert a new daily price record as row 14, pushing the
# existing rows 14-77 down to 15-78 (dimension grows from A1:R77 to A1:R78).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 14 (shifts rows 14..77 -> 15..78).
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with the new record.
$ws.Cells.Item(14, 1).Value  = 4
$ws.Cells.Item(14, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(14, 3).Value  = "Los Lagos"
$ws.Cells.Item(14, 4).Value  = 44575
$ws.Cells.Item(14, 5).Value  = 10
$ws.Cells.Item(14, 6).Value  = 100112026
$ws.Cells.Item(14, 7).Value  = "Haba"
$ws.Cells.Item(14, 8).Value  = "Sin especificar"
$ws.Cells.Item(14, 9).Value  = "Primera"
$ws.Cells.Item(14, 10).Value = 80
$ws.Cells.Item(14, 11).Value = 24000
$ws.Cells.Item(14, 12).Value = 24000
$ws.Cells.Item(14, 13).Value = 24000
$ws.Cells.Item(14, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(14, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(14, 16).Value = 960
$ws.Cells.Item(14, 17).Value = 25
$ws.Cells.Item(14, 18).Value = "Hortaliza"
